$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 83

$ws.Cells.Item($newRow, 1).Value = 46032
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

$ws.Cells.Item($newRow, 2).Value = 185
$ws.Cells.Item($newRow, 3).Value = 198
$ws.Cells.Item($newRow, 4).Value = 189
